$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.247.28"
$ws.Range("D3").Value = "1.864.19"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2868"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +10.34%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").Value = "1.871.10"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.192"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6834"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "277.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "30.247.90"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.83%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007348"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "2.119.03"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.361"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.200"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.250"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.958"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09848"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.383"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04748"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +4.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7051"
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("E39").Value = "  +4.63%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.289"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.959"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8537"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.217"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "949.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05646"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
